# Updates cryptos list values (price and volume columns) per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '57.155.44' },
    @{ Cell = 'E2'; Value = '  +0.30%  ' },
    @{ Cell = 'D3'; Value = '3.061.73' },
    @{ Cell = 'E3'; Value = '  +1.10%  ' },
    @{ Cell = 'E4'; Value = '  -0.02%  ' },
    @{ Cell = 'D5'; Value = '515.65' },
    @{ Cell = 'E5'; Value = '  +0.56%  ' },
    @{ Cell = 'D6'; Value = '140.34' },
    @{ Cell = 'E6'; Value = '  -0.16%  ' },
    @{ Cell = 'D7'; Value = '1.00' },
    @{ Cell = 'E7'; Value = '  -0.07%  ' },
    @{ Cell = 'D8'; Value = '0.432' },
    @{ Cell = 'E8'; Value = '  -1.55%  ' },
    @{ Cell = 'D9'; Value = '7.24' },
    @{ Cell = 'E9'; Value = '  +0.73%  ' },
    @{ Cell = 'E10'; Value = '  -0.99%  ' },
    @{ Cell = 'D11'; Value = '0.369' },
    @{ Cell = 'E11'; Value = '  -1.59%  ' },
    @{ Cell = 'D12'; Value = '3.596.51' },
    @{ Cell = 'E12'; Value = '  +0.80%  ' },
    @{ Cell = 'E13'; Value = '  +2.65%  ' },
    @{ Cell = 'D14'; Value = '25.42' },
    @{ Cell = 'E14'; Value = '  -4.57%  ' },
    @{ Cell = 'D15'; Value = '0.0000162' },
    @{ Cell = 'E15'; Value = '  -1.81%  ' },
    @{ Cell = 'D16'; Value = '57.256.75' },
    @{ Cell = 'E16'; Value = '  +0.53%  ' },
    @{ Cell = 'D17'; Value = '3.073.22' },
    @{ Cell = 'E17'; Value = '  +1.20%  ' },
    @{ Cell = 'D18'; Value = '6.05' },
    @{ Cell = 'E18'; Value = '  -0.70%  ' },
    @{ Cell = 'D19'; Value = '12.96' },
    @{ Cell = 'E19'; Value = '  -2.28%  ' },
    @{ Cell = 'D20'; Value = '8.04' },
    @{ Cell = 'E20'; Value = '  +0.34%  ' },
    @{ Cell = 'D21'; Value = '332.34' },
    @{ Cell = 'E21'; Value = '  +0.12%  ' },
    @{ Cell = 'E22'; Value = '  -0.20%  ' },
    @{ Cell = 'D23'; Value = '0.498' },
    @{ Cell = 'E23'; Value = '  -0.91%  ' },
    @{ Cell = 'D24'; Value = '65.49' },
    @{ Cell = 'E24'; Value = '  +0.26%  ' },
    @{ Cell = 'E25'; Value = '  +3.27%  ' },
    @{ Cell = 'E26'; Value = '  -0.54%  ' },
    @{ Cell = 'D27'; Value = '0.0₃0903' },
    @{ Cell = 'E27'; Value = '  +0.28%  ' },
    @{ Cell = 'D28'; Value = '6.32' },
    @{ Cell = 'E28'; Value = '  -5.25%  ' },
    @{ Cell = 'D29'; Value = '7.12' },
    @{ Cell = 'E29'; Value = '  -0.07%  ' },
    @{ Cell = 'E30'; Value = '  +0.66%  ' },
    @{ Cell = 'D31'; Value = '20.75' },
    @{ Cell = 'E31'; Value = '  +0.79%  ' },
    @{ Cell = 'D32'; Value = '1.16' },
    @{ Cell = 'E32'; Value = '  -3.03%  ' },
    @{ Cell = 'D33'; Value = '154.77' },
    @{ Cell = 'E33'; Value = '  +1.00%  ' },
    @{ Cell = 'B34'; Value = 'EnergySwap' },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D34'; Value = '27.28' },
    @{ Cell = 'E34'; Value = '  +8.61%  ' },
    @{ Cell = 'B35'; Value = 'NEARProtocol' },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Cell = 'D35'; Value = '4.45' },
    @{ Cell = 'E35'; Value = '  -4.70%  ' },
    @{ Cell = 'D36'; Value = '5.82' },
    @{ Cell = 'E36'; Value = '  -1.52%  ' },
    @{ Cell = 'D37'; Value = '1.27' },
    @{ Cell = 'E37'; Value = '  +0.30%  ' },
    @{ Cell = 'D38'; Value = '0.0669' },
    @{ Cell = 'E38'; Value = '  -0.15%  ' },
    @{ Cell = 'D39'; Value = '3.109.16' },
    @{ Cell = 'E39'; Value = '  +1.36%  ' },
    @{ Cell = 'D40'; Value = '36.85' },
    @{ Cell = 'E40'; Value = '  -0.31%  ' },
    @{ Cell = 'D41'; Value = '3.87' },
    @{ Cell = 'E41'; Value = '  +0.11%  ' },
    @{ Cell = 'D42'; Value = '0.999' },
    @{ Cell = 'E42'; Value = '  -0.03%  ' },
    @{ Cell = 'D43'; Value = '0.657' },
    @{ Cell = 'E43'; Value = '  -0.66%  ' },
    @{ Cell = 'D44'; Value = '2.264.78' },
    @{ Cell = 'E44'; Value = '  +3.06%  ' },
    @{ Cell = 'D45'; Value = '0.0258' },
    @{ Cell = 'E45'; Value = '  +6.75%  ' },
    @{ Cell = 'E46'; Value = '  -1.76%  ' },
    @{ Cell = 'D47'; Value = '0.926' },
    @{ Cell = 'E47'; Value = '  -2.90%  ' },
    @{ Cell = 'D48'; Value = '5.85' },
    @{ Cell = 'E48'; Value = '  -2.22%  ' },
    @{ Cell = 'D49'; Value = '19.76' },
    @{ Cell = 'E49'; Value = '  -2.20%  ' },
    @{ Cell = 'D50'; Value = '0.0870' },
    @{ Cell = 'E50'; Value = '  +0.53%  ' },
    @{ Cell = 'D51'; Value = '249.28' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Preserve original text representation (e.g. "57.155.44", "1.00", "0.0000162")
    # instead of letting Excel auto-convert numeric-looking strings to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
